{"js": "// The document contains a paragraph whose text is the M2Doc user-content\n// field marker:  {m:userdoc 'zone1'}\n// It currently lives in two runs: \"{m\" and \":userdoc 'zone1'}\".\n// The expected (migrated) result splits that same text across four runs:\n//   \"{\" | \"m\" | \":userdoc 'zone1'\" | \"}\"\n// (with the final run keeping xml:space=\"preserve\"), and the run/paragraph\n// rsid bookkeeping attributes are dropped, matching the parser's new\n// TokenIteratorFieldRewriterSplit behaviour.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph holding the \"{m:userdoc 'zone1'}\" field text.\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"{m\") !== -1 && p.text.indexOf(\"userdoc\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  const range = target.getRange();\n\n  // Rebuild the paragraph content as four distinct runs by inserting an\n  // OOXML fragment in place of the current range - this is what forces the\n  // single-run-per-text-chunk boundaries the diff expects.\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>m</w:t></w:r>' +\n    '<w:r><w:t>:userdoc \\'zone1\\'</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n\n  range.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The document contains a paragraph whose text is the M2Doc user-content\n# field marker:  {m:userdoc 'zone1'}\n# It currently lives in two runs: \"{m\" and \":userdoc 'zone1'}\".\n# The expected (migrated) result splits that same text across four runs:\n#   \"{\" | \"m\" | \":userdoc 'zone1'\" | \"}\"\n# (with the final run keeping xml:space=\"preserve\"), and the run/paragraph\n# rsid bookkeeping attributes are dropped, matching the parser's new\n# TokenIteratorFieldRewriterSplit behaviour.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*{m*\" -and $t -like \"*userdoc*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $r = $target.Range\n\n    $xml = \"<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\" +\n           \"<w:r><w:t>{</w:t></w:r>\" +\n           \"<w:r><w:t>m</w:t></w:r>\" +\n           \"<w:r><w:t>:userdoc 'zone1'</w:t></w:r>\" +\n           \"<w:r><w:t xml:space='preserve'>}</w:t></w:r>\" +\n           \"</w:p>\"\n\n    # InsertXML replaces the exact range's contents with the supplied runs,\n    # giving precise control over the run boundaries (unlike Range.Text,\n    # which merges everything back into a single run).\n    $r.InsertXML($xml)\n}\n"}
